# HEVC FEI 16x16 MVP skip-mode deck:
#   - Fix Figure#4 zigzag order: swap the "5"/"6" labels on the
#     flow-chart process boxes (picture replaced with correct zigzag
#     order).
#   - Refresh the cached "today" date shown by the auto-updating
#     date-time footer field on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Figure#4 zig-zag numbering fix (single slide in this deck).
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }

    $text = $shp.TextFrame.TextRange.Text

    if ($shp.Name -eq "Flowchart: Process 82" -and $text -eq "6") {
        # Was mislabeled "6" - should read "5".
        $shp.TextFrame.TextRange.Text = "5"
    }
    elseif ($shp.Name -eq "Flowchart: Process 86" -and $text -eq "5") {
        # Was mislabeled "5" - should read "6".
        $shp.TextFrame.TextRange.Text = "6"
    }
}

# ---------------------------------------------------------------------
# 2. Re-cache the datetimeFigureOut placeholder text (master + every
#    slide layout) to the current save date.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = "4/18/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
